$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.520.96"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.048.67"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.37"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.658"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.41"
$ws.Range("E8").Value = "  -7.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.51"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0748"
$ws.Range("E11").Value = "  -4.25%  "
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.952"
$ws.Range("E13").Value = "  +7.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.71"
$ws.Range("E14").Value = "  -4.39%  "
$ws.Range("D15").Value = "2.348.55"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.43"
$ws.Range("E16").Value = "  -4.63%  "
$ws.Range("D17").Value = "2.040.44"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "36.390.71"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.03"
$ws.Range("E19").Value = "  -6.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.70"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").Value = "0.0₃0856"
$ws.Range("E21").Value = "  -4.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.12"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("E23").Value = "  -4.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("E26").Value = "  +3.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.92"
$ws.Range("E27").Value = "  -2.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.15"
$ws.Range("E28").Value = "  -9.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.97"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").Value = "  +9.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.04"
$ws.Range("E32").Value = "  -9.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.44"
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0591"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0860"
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  -5.82%  "
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.23"
$ws.Range("E40").Value = "  -7.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("E41").Value = "  -5.38%  "
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("E43").Value = "  -5.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.05"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0904"
$ws.Range("E45").Value = "  -5.23%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.408.32"
$ws.Range("E46").Value = "  +7.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.78"
$ws.Range("E47").Value = "  -6.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.42"
$ws.Range("E48").Value = "  +9.03%  "
$ws.Range("E50").Value = "  -5.04%  "
$ws.Range("D51").Value = "2.230.68"
$ws.Range("E51").Value = "  -0.63%  "
